# Apply the edits described by the diff:
#  - Column G ("target") data values change from "proton" to "p" for all
#    data rows (rows 2-13 on Sheet1).
#  - The header row (row 1) is made bold.
#  - The sheet selection is changed to G2:G13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the "target" column (G) values from "proton" to "p" ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 13 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    if ($cell.Value2 -eq "proton") {
        $cell.Value = "p"
    }
}

# --- Make the header row bold ---
$ws.Range("A1:K1").Font.Bold = $true

# --- Update the selected range shown in the sheet view ---
$ws.Activate()
$ws.Range("G2:G13").Select()
